$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a Price (column D) cell to a literal text value without Excel
# re-interpreting decimal-looking strings (e.g. "1.000", "39.74") as numbers.
# coinranking values are pre-formatted display strings, always stored as text
# in the sheet (t="inlineStr"/shared-string), so we force text via NumberFormat
# "@" for the assignment, then restore the cell style (no explicit number format
# was set on these cells originally).
function Set-PriceText($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- Simple per-row Price (D) / Volume(1h) (E) updates ---
Set-PriceText $ws.Range("D2") "19.950.76"
$ws.Range("E2").Value = "  -5.44%  "
Set-PriceText $ws.Range("D3") "1.413.96"
$ws.Range("E3").Value = "  -6.31%  "
Set-PriceText $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  -0.64%  "
Set-PriceText $ws.Range("D5") "1.000"
$ws.Range("E5").Value = "  -0.66%  "
Set-PriceText $ws.Range("D6") "276.06"
$ws.Range("E6").Value = "  -3.18%  "
Set-PriceText $ws.Range("D7") "0.3668"
$ws.Range("E7").Value = "  -5.19%  "
Set-PriceText $ws.Range("D8") "0.3101"
$ws.Range("E8").Value = "  -0.88%  "
Set-PriceText $ws.Range("D9") "39.73"
$ws.Range("E9").Value = "  -5.84%  "
$ws.Range("E10").Value = "  -0.53%  "
Set-PriceText $ws.Range("D11") "0.06528"
$ws.Range("E11").Value = "  -6.57%  "
Set-PriceText $ws.Range("D12") "1.001"
$ws.Range("E12").Value = "  -0.81%  "
Set-PriceText $ws.Range("D13") "5.485"
$ws.Range("E13").Value = "  -2.34%  "
Set-PriceText $ws.Range("D14") "17.61"
$ws.Range("E14").Value = "  -1.45%  "
Set-PriceText $ws.Range("D15") "6.192"
$ws.Range("E15").Value = "  -2.90%  "
Set-PriceText $ws.Range("D16") "1.412.83"
$ws.Range("E16").Value = "  -6.94%  "
Set-PriceText $ws.Range("D17") "0.00001017"
$ws.Range("E17").Value = "  -4.89%  "
Set-PriceText $ws.Range("D18") "0.05662"
$ws.Range("E18").Value = "  -13.92%  "
Set-PriceText $ws.Range("D19") "1.001"
$ws.Range("E19").Value = "  -0.63%  "
Set-PriceText $ws.Range("D20") "71.18"
$ws.Range("E20").Value = "  -13.25%  "
Set-PriceText $ws.Range("D21") "5.619"
$ws.Range("E21").Value = "  -6.43%  "
Set-PriceText $ws.Range("D22") "14.74"
$ws.Range("E22").Value = "  -2.60%  "
Set-PriceText $ws.Range("D23") "10.92"
$ws.Range("E23").Value = "  +0.96%  "
Set-PriceText $ws.Range("D24") "2.237"
$ws.Range("E24").Value = "  -4.97%  "
Set-PriceText $ws.Range("D25") "19.958.54"
$ws.Range("E25").Value = "  -5.39%  "
Set-PriceText $ws.Range("D26") "2.263"
$ws.Range("E26").Value = "  -3.48%  "
Set-PriceText $ws.Range("D27") "132.98"
$ws.Range("E27").Value = "  -9.97%  "
Set-PriceText $ws.Range("D28") "17.29"
$ws.Range("E28").Value = "  -4.06%  "
Set-PriceText $ws.Range("D29") "1.569.75"
$ws.Range("E29").Value = "  -6.98%  "
Set-PriceText $ws.Range("D30") "109.63"
$ws.Range("E30").Value = "  -4.40%  "
Set-PriceText $ws.Range("D31") "3.899"
$ws.Range("E31").Value = "  -18.67%  "
Set-PriceText $ws.Range("D32") "5.280"
$ws.Range("E32").Value = "  -11.27%  "
Set-PriceText $ws.Range("D33") "0.8206"
$ws.Range("E33").Value = "  -13.89%  "
Set-PriceText $ws.Range("D34") "0.07695"
$ws.Range("E34").Value = "  -3.25%  "
Set-PriceText $ws.Range("D35") "1.485"
$ws.Range("E35").Value = "  +0.70%  "
Set-PriceText $ws.Range("D36") "8.329"
$ws.Range("E36").Value = "  -0.73%  "
Set-PriceText $ws.Range("D37") "4.916"
$ws.Range("E37").Value = "  -3.18%  "
Set-PriceText $ws.Range("D38") "0.05788"
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("E39").Value = "  -0.60%  "
Set-PriceText $ws.Range("D40") "0.02060"
$ws.Range("E40").Value = "  -3.39%  "
Set-PriceText $ws.Range("D41") "10.46"
$ws.Range("E41").Value = "  -7.31%  "
Set-PriceText $ws.Range("D42") "0.1883"
$ws.Range("E42").Value = "  -4.47%  "
Set-PriceText $ws.Range("D43") "1.095"
$ws.Range("E43").Value = "  -5.31%  "
Set-PriceText $ws.Range("D47") "0.5187"
$ws.Range("E47").Value = "  -4.30%  "
Set-PriceText $ws.Range("D48") "115.77"
$ws.Range("E48").Value = "  +1.84%  "
Set-PriceText $ws.Range("D49") "1.769"
$ws.Range("E49").Value = "  -4.16%  "
Set-PriceText $ws.Range("D50") "1.033"
$ws.Range("E50").Value = "  -9.17%  "
Set-PriceText $ws.Range("D51") "1.000"
$ws.Range("E51").Value = "  -0.78%  "

# --- Rows 44-46: coin order reshuffled (EnergySwap, TheSandbox, PancakeSwap) ---
# with refreshed Coin / Link / Price / Volume(1h) values
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-PriceText $ws.Range("D44") "12.45"
$ws.Range("E44").Value = "  -3.59%  "

$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-PriceText $ws.Range("D45") "0.5313"
$ws.Range("E45").Value = "  -5.37%  "

$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-PriceText $ws.Range("D46") "3.537"
$ws.Range("E46").Value = "  -4.02%  "

